$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style from H1 (bold, bordered, centered) onto the new
# header cells I1 and J1 before setting their text, so the format (s="1")
# travels along with the copy.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Header labels for the two new columns
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# New data values for columns I (I0) and J (IF)
$ws.Range("I2").Value = 4
$ws.Range("J2").Value = 5

$ws.Range("I3").Value = 5
$ws.Range("J3").Value = 7

$ws.Range("I4").Value = 5
$ws.Range("J4").Value = 7

$ws.Range("I5").Value = 6
$ws.Range("J5").Value = 8

$ws.Range("I6").Value = 1
$ws.Range("J6").Value = 3

$ws.Range("I7").Value = 3
$ws.Range("J7").Value = 4
